$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.3301315503884687
$ws.Range("B2").Value = 0.2955902794325816
$ws.Range("C2").Value = 0.3630767971809157
$ws.Range("D2").Value = 0.2725162596895327
$ws.Range("E2").Value = 0.4066187414426702
$ws.Range("K2").Value = 2.31092085271928
$ws.Range("L2").Value = 2.069131956028071
$ws.Range("M2").Value = 2.54153758026641
$ws.Range("N2").Value = 1.907613817826729
$ws.Range("O2").Value = 2.846331190098691
$ws.Range("P2").Value = 64.44223
$ws.Range("Q2").Value = 46.30304518419486
$ws.Range("R2").Value = 87.33486795507919
$ws.Range("S2").Value = 53.5777423435626
$ws.Range("T2").Value = 74.5312578753298
$ws.Range("U2").Value = 0.2031548274478529
$ws.Range("V2").Value = 0.164662504147416
$ws.Range("W2").Value = 0.241878409505886
$ws.Range("X2").Value = 0.1361939362162032
$ws.Range("Y2").Value = 0.2795835804766086
$ws.Range("Z2").Value = 0.8573167908766585
$ws.Range("AA2").Value = 0.8113466238081585
$ws.Range("AB2").Value = 0.8980344071975621
$ws.Range("AC2").Value = 0.7709565353244547
$ws.Range("AD2").Value = 0.9298609587351283
$ws.Range("F3").Value = 7.001534600731914
$ws.Range("G3").Value = 6.220797888576329
$ws.Range("H3").Value = 7.754317288312638
$ws.Range("I3").Value = 5.698442768810907
$ws.Range("J3").Value = 8.7717687647198
$ws.Range("K3").Value = 2.310506418241533
$ws.Range("L3").Value = 2.052863303230189
$ws.Range("M3").Value = 2.558924705143171
$ws.Range("N3").Value = 1.880486113707599
$ws.Range("O3").Value = 2.894683692357535
$ws.Range("P3").Value = 63.51508
$ws.Range("Q3").Value = 56.40103009049271
$ws.Range("R3").Value = 72.22926416319957
$ws.Range("S3").Value = 59.40179431101897
$ws.Range("T3").Value = 67.77768657085397
$ws.Range("U3").Value = 0.2028513684196727
$ws.Range("V3").Value = 0.1617836624909226
$ws.Range("W3").Value = 0.2444918034730678
$ws.Range("X3").Value = 0.131574774878728
$ws.Range("Y3").Value = 0.2857563357823072
$ws.Range("Z3").Value = 0.856172181610088
$ws.Range("AA3").Value = 0.8067508856446146
$ws.Range("AB3").Value = 0.8998751425036341
$ws.Range("AC3").Value = 0.7631267830679144
$ws.Range("AD3").Value = 0.9340652383046261
$ws.Range("A4").Value = 0.3299527873842167
$ws.Range("B4").Value = 0.2955532292028652
$ws.Range("C4").Value = 0.3629908133231589
$ws.Range("D4").Value = 0.2719585148209576
$ws.Range("E4").Value = 0.4067924091754813
$ws.Range("F4").Value = 7.003102622518629
$ws.Range("G4").Value = 6.217063611490602
$ws.Range("H4").Value = 7.758896798529163
$ws.Range("I4").Value = 5.700984323911813
$ws.Range("J4").Value = 8.774824090174556
$ws.Range("K4").Value = 2.310581376869882
$ws.Range("L4").Value = 1.960099274045088
$ws.Range("M4").Value = 2.647115016992029
$ws.Range("N4").Value = 1.725673933688124
$ws.Range("O4").Value = 3.109328671323687
$ws.Range("P4").Value = 65.36592
$ws.Range("Q4").Value = 45.65299091608324
$ws.Range("R4").Value = 94.4306790924532
$ws.Range("S4").Value = 53.76457611257006
$ws.Range("T4").Value = 76.43869195732356
$ws.Range("U4").Value = 0.2013647291966062
$ws.Range("V4").Value = 0.14530336231264
$ws.Range("W4").Value = 0.2575543796331879
$ws.Range("X4").Value = 0.1034315077128284
$ws.Range("Y4").Value = 0.3121750978671887
$ws.Range("Z4").Value = 0.8490013627832288
$ws.Range("AA4").Value = 0.7776481050978237
$ws.Range("AB4").Value = 0.9080409396528987
$ws.Range("AC4").Value = 0.7066943562492456
$ws.Range("AD4").Value = 0.9517019068147211
